$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header "MemoryId" (column D) to "MemorycardId"
$ws.Range("D1").Value = "MemorycardId"

# Move selection to D1 to match the saved view state
$ws.Range("D1").Select()
